# Update countries & provincias Spain
# Applies the data refresh captured in the commit: a handful of countries'
# case counts were updated, which shuffled their rank (rows are kept sorted
# by "Casos totales" descending), plus the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Row 15 - Iran: refreshed stats, no rank change
Set-Row 15 "Iran" 395488 2063 340842 31848 0 129 22798

# Rows 36-37 - Rumania overtakes Panama
Set-Row 36 "Rumania" 99684 1380 41010 54609 0 47 4065
Set-Row 37 "Panama" 99042 0 71419 25507 0 0 2116

# Rows 60-62 - Suiza overtakes Armenia and Ghana
Set-Row 60 "Suiza" 45711 405 38100 5592 0 0 2019
Set-Row 61 "Armenia" 45326 174 41233 3187 0 1 906
Set-Row 62 "Ghana" 45313 0 44188 842 0 0 283

# Row 85 - Madagascar: refreshed stats, no rank change
Set-Row 85 "Madagascar" 15624 104 14295 1121 0 2 208

# Row 146 - Malta: refreshed stats, no rank change
Set-Row 146 "Malta" 2204 42 1803 387 0 0 14

# Row 183 - Gibraltar: refreshed stats, no rank change
Set-Row 183 "Gibraltar" 323 1 287 36 0 0 0

# Rows 206-209 - Bonaire, San Eustaquio y Saba overtakes Granada and Laos;
# Dominica keeps its spot but gets refreshed stats
Set-Row 206 "Bonaire, San Eustaquio y Saba" 25 4 7 18 0 0 0
Set-Row 207 "Dominica" 24 2 18 6 0 0 0
Set-Row 208 "Granada" 24 0 24 0 0 0 0
Set-Row 209 "Laos" 22 0 21 1 0 0 0

# Rows 214-215 - Montserrat overtakes Islas Malvinas
Set-Row 214 "Montserrat" 13 0 12 0 0 0 1
Set-Row 215 "Islas Malvinas" 13 0 13 0 0 0 0

# Rows 217-218 - San Pedro y Miquelon overtakes Sahara Occidental
Set-Row 217 "San Pedro y Miquelon" 10 1 5 5 0 0 0
Set-Row 218 "Sahara Occidental" 10 0 8 1 0 0 1

# Update the "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 13:34"
